# "lectura de datos realizado" - fill in the data that was read/captured
# for the last two students on the "Sheet" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Nogueda Hernadez Hugo David) was missing its "semestre" value.
$ws.Range("E9").Value = 6

# Row 10 is a brand new record: same person as row 8 (Lopez Mendoza
# Tania Guadalupe), semestre 6.
$ws.Range("A10").Value = "Lopez"
$ws.Range("B10").Value = "Mendoza"
# Column B carries an inherited column style; reset this cell back to
# the workbook's default "Normal" style so it doesn't pick that up.
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "Tania"
$ws.Range("D10").Value = "Guadalupe"
$ws.Range("E10").Value = 6
